$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.013.11"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").Value = "3.486.02"
$ws.Range("E3").Value = "  -2.23%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "585.15"
$ws.Range("E5").Value = "  +5.91%  "

$ws.Range("D6").Value = "177.96"
$ws.Range("E6").Value = "  -6.01%  "

$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +3.81%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.636"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("E10").Value = "  +5.22%  "

$ws.Range("D11").Value = "56.06"
$ws.Range("E11").Value = "  +2.08%  "

$ws.Range("D12").Value = "0.0000278"
$ws.Range("E12").Value = "  +2.57%  "

$ws.Range("D13").Value = "9.27"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").Value = "4.043.40"
$ws.Range("E14").Value = "  -2.22%  "

$ws.Range("D15").Value = "3.484.84"
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "18.29"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").Value = "12.08"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "65.971.33"
$ws.Range("E19").Value = "  -1.66%  "

$ws.Range("E20").Value = "  +1.43%  "

$ws.Range("D21").Value = "413.01"
$ws.Range("E21").Value = "  -4.85%  "

$ws.Range("D22").Value = "4.33"
$ws.Range("E22").Value = "  +10.51%  "

$ws.Range("E23").Value = "  +6.18%  "

$ws.Range("D24").Value = "84.76"
$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("D25").Value = "13.45"
$ws.Range("E25").Value = "  +10.84%  "

$ws.Range("D26").Value = "11.06"
$ws.Range("E26").Value = "  -0.47%  "

$ws.Range("D27").Value = "2.87"
$ws.Range("E27").Value = "  -1.34%  "

$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("D29").Value = "9.21"
$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("D30").Value = "30.28"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").Value = "6.68"
$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("D32").Value = "11.75"
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").Value = "594.52"
$ws.Range("E33").Value = "  -8.46%  "

$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").Value = "61.01"
$ws.Range("E35").Value = "  +2.02%  "

$ws.Range("E36").Value = "  +0.65%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").Value = "0.0₃0796"
$ws.Range("E38").Value = "  -3.79%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "36.94"
$ws.Range("E39").Value = "  -4.57%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "3.59"
$ws.Range("E40").Value = "  +6.40%  "

$ws.Range("D41").Value = "0.385"
$ws.Range("E41").Value = "  -1.83%  "

$ws.Range("D42").Value = "3.223.66"
$ws.Range("E42").Value = "  +6.19%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "2.98"
$ws.Range("E44").Value = "  +3.23%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  -5.37%  "

$ws.Range("D47").Value = "0.0420"
$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("D49").Value = "2.66"
$ws.Range("E49").Value = "  -4.65%  "

$ws.Range("D50").Value = "8.60"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("D51").Value = "139.98"
$ws.Range("E51").Value = "  -2.64%  "
